# The MAPE columns in this workbook were stored as raw fractions (e.g. 0.15)
# but should be expressed as percentages (e.g. 15.0). Multiply the relevant
# MAPE cells by 100 on both sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "metrics_food": columns F (MAPE_Q1), G (MAPE_Q2), H (MAPE_Q3) ---
# Header is on row 1, data runs from row 2 through row 75.
$wsFood = $wb.Worksheets.Item("metrics_food")

for ($r = 2; $r -le 75; $r++) {
    foreach ($col in @("F", "G", "H")) {
        $cell = $wsFood.Range("$col$r")
        $cell.Value = $cell.Value2 * 100
    }
}

# --- Sheet "metrics_agg": column C (MAPE), rows 2 through 4 ---
$wsAgg = $wb.Worksheets.Item("metrics_agg")

for ($r = 2; $r -le 4; $r++) {
    $cell = $wsAgg.Range("C$r")
    $cell.Value = $cell.Value2 * 100
}
